$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price observation inserted above the existing row 41 ("Hortaliza,
# Bruselas (repollito)" series). Excel shifts every row from 41..84 down to
# 42..85, preserving all of their original values - this single insert
# reproduces the whole cascade of diffs shown for rows 41-85.
$ws.Rows.Item(41).Insert()

# Populate the freshly inserted row 41 with the new record's data.
$ws.Range("A41").Value = 9
$ws.Range("B41").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 45079
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 100112035
$ws.Range("G41").Value = "Bruselas (repollito)"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 52
$ws.Range("K41").Value = 19000
$ws.Range("L41").Value = 21000
$ws.Range("M41").Value = 20000
$ws.Range("N41").Value = "$/malla 15 kilos"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 1333
$ws.Range("Q41").Value = 15
$ws.Range("R41").Value = "Hortaliza"
